$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently displays the rule-id text "R40" (a shared string).
# The saved/updated rule now uses "1" as its id, so update the cell's
# text accordingly. The leading apostrophe forces Excel to store this
# as literal text rather than re-interpreting the numeric-looking
# string as a number, matching the original cell's string type.
$ws.Range("B11").Value = "'1"
